# Generate Report for Handback
# Applies the handback results to the localization-status workbook:
#  - updates the "Status" text from "Ready for handoff" to
#    "Handed back: in sync with en-US" on all three sheets
#  - fills in the "Latest Handback File" / "Latest Handback DateTime"
#    columns (and hyperlinks) for the zh-cn and de-de sheets
#  - refreshes the previously-unset handback datetime placeholder
#  - widens a few columns so the new content is readable

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: status text for both locales, both files ---
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText

# --- Overview sheet: widen the zh-cn / de-de status columns ---
$ws1.Columns.Item(5).ColumnWidth = 29.15
$ws1.Columns.Item(6).ColumnWidth = 29.15

# --- zh-cn sheet ---
$ws2.Range("C2").Value = $statusText
$ws2.Range("C3").Value = $statusText

$ws2.Range("I2").Value = "07f6527b-655e-46b8-9151-70c1b8059072.md"
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1111b4de758a143abfd6751fb2be7904f5d2b4e/e2e/07f6527b-655e-46b8-9151-70c1b8059072.md", "", "", "07f6527b-655e-46b8-9151-70c1b8059072.md") | Out-Null
$ws2.Range("I2").Font.Underline = $true
$ws2.Range("I2").Font.Color = 15570276
$ws2.Range("J2").Value = "07f6527b-655e-46b8-9151-70c1b8059072.10e1a39fd77831289fdaa5a8da6546f6c310dc80.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-06 21:05:49"

$ws2.Range("I3").Value = "6135c9e7-7e16-4819-8d49-20d1572f7e07.md"
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1111b4de758a143abfd6751fb2be7904f5d2b4e/e2e/6135c9e7-7e16-4819-8d49-20d1572f7e07.md", "", "", "6135c9e7-7e16-4819-8d49-20d1572f7e07.md") | Out-Null
$ws2.Range("I3").Font.Underline = $true
$ws2.Range("I3").Font.Color = 15570276
$ws2.Range("J3").Value = "6135c9e7-7e16-4819-8d49-20d1572f7e07.2498bf6843f1b44ef93893d0cf0750c8d8fd24cb.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-09-06 21:05:49"

# --- zh-cn sheet: widen columns ---
$ws2.Columns.Item(3).ColumnWidth = 29.15
$ws2.Columns.Item(9).ColumnWidth = 39.15
$ws2.Columns.Item(10).ColumnWidth = 39.15

# --- de-de sheet ---
$ws3.Range("C2").Value = $statusText
$ws3.Range("C3").Value = $statusText

$ws3.Range("I2").Value = "07f6527b-655e-46b8-9151-70c1b8059072.md"
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1111b4de758a143abfd6751fb2be7904f5d2b4e/e2e/07f6527b-655e-46b8-9151-70c1b8059072.md", "", "", "07f6527b-655e-46b8-9151-70c1b8059072.md") | Out-Null
$ws3.Range("I2").Font.Underline = $true
$ws3.Range("I2").Font.Color = 15570276
$ws3.Range("J2").Value = "07f6527b-655e-46b8-9151-70c1b8059072.10e1a39fd77831289fdaa5a8da6546f6c310dc80.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-06 21:05:57"

$ws3.Range("I3").Value = "6135c9e7-7e16-4819-8d49-20d1572f7e07.md"
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1111b4de758a143abfd6751fb2be7904f5d2b4e/e2e/6135c9e7-7e16-4819-8d49-20d1572f7e07.md", "", "", "6135c9e7-7e16-4819-8d49-20d1572f7e07.md") | Out-Null
$ws3.Range("I3").Font.Underline = $true
$ws3.Range("I3").Font.Color = 15570276
$ws3.Range("J3").Value = "6135c9e7-7e16-4819-8d49-20d1572f7e07.2498bf6843f1b44ef93893d0cf0750c8d8fd24cb.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-06 21:05:57"

# --- de-de sheet: widen columns ---
$ws3.Columns.Item(3).ColumnWidth = 29.15
$ws3.Columns.Item(9).ColumnWidth = 39.15
$ws3.Columns.Item(10).ColumnWidth = 39.15
